# TODO.xlsx update:
# - Task "Verifizierungs-E-Mail nach Registierung" (row 5) is now done:
#   status changed from "in Arbeit" to "done" (reuse formatting from an
#   existing "done" cell so fill/font match the other done rows).
# - Selection/active cell left on D12 when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing "done" status cell (B2) onto B5,
# then set its text to "done".
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value = "done"

# Clear the clipboard marching ants / copy mode.
$excel.CutCopyMode = 0

# Move the selection to D12, matching where the cursor was left on save.
$ws.Range("D12").Select() | Out-Null
